$p = $ppt.ActivePresentation

# Update the literal date text shown in the footer's "Date Placeholder"
# on the slide master and on every slide layout (8/10/2023 -> 8/16/2023).
$p.SlideMaster.HeadersFooters.DateAndTime.Text = "8/16/2023"

for ($i = 1; $i -le $p.SlideMaster.CustomLayouts.Count; $i++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($i)
    $layout.HeadersFooters.DateAndTime.Text = "8/16/2023"
}
